$wb = $excel.ActiveWorkbook
$fullSheet2018 = $wb.Worksheets.Item("2018-Full")
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $fullSheet2018)
$new.Name = "2019-Full"

$new.Range("A1").Value = "State"
$new.Range("A1").Font.Name = "Verdana"
$new.Range("A1").Borders.LineStyle = 1
$new.Range("A1").VerticalAlignment = -4108

$new.Range("B2").Value = 124330
$new.Range("B2").NumberFormat = "#,##0 ;(#,##0)"
$new.Range("B2").Borders.LineStyle = 1

Write-Output "done"
